$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the "Total" summary row (currently row 21) down to row 23, two rows
#    lower, keeping its formatting (grey fill / bold / borders) and its
#    formulas (which keep referencing D2:D20 / F2:F20 verbatim - they are not
#    re-targeted to the newly inserted rows below).
$ws.Range("A21:F21").Copy()
$ws.Range("A23:F23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Total"
$ws.Range("D23").Formula = "=SUM(D2:D20)"
$ws.Range("F23").Formula = "=SUM(F2:F20)"

# 2. Add two new timesheet entries as rows 21 and 22, copying the formatting
#    of the last existing data row (row 20) so the date/time/number formats
#    and borders match the rest of the table.
$ws.Range("A20:F20").Copy()
$ws.Range("A21:F21").PasteSpecial(-4122)
$ws.Range("A22:F22").PasteSpecial(-4122)

$ws.Range("A21").Value = 45311
$ws.Range("B21").Value = 0.791666666666667
$ws.Range("C21").Value = 0.958333333333333
$ws.Range("D21").Formula = "=(C21<B21)+C21-B21"
$ws.Range("E21").Value = 13.5
$ws.Range("F21").Formula = "=(D21*24)*E21"
$ws.Range("F21").NumberFormat = "General"

$ws.Range("A22").Value = 45312
$ws.Range("B22").Value = 0.5
$ws.Range("C22").Value = 0.583333333333333
$ws.Range("D22").Formula = "=(C22<B22)+C22-B22"
$ws.Range("E22").Value = 13.5
$ws.Range("F22").Formula = "=(D22*24)*E22"
$ws.Range("F22").NumberFormat = "General"

# 3. Bump the hourly rate from 10 to 13.5 for every existing timesheet entry
#    (rows 2-20); the Bill column (F) recalculates automatically since it is
#    a formula.
for ($r = 2; $r -le 20; $r++) {
    $ws.Range("E$r").Value = 13.5
}

# 4. Match the author's final selection (cell D22) instead of the old F22.
$ws.Range("D22").Select() | Out-Null
